$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsZh.Range("D3").Value = "2016-02-16 14:56:34"
$wsZh.Range("G3").Value = "2016-02-16 14:57:31"

$wsDe.Range("D3").Value = "2016-02-16 14:56:51"
$wsDe.Range("G3").Value = "2016-02-16 14:58:02"
